$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54, shifting existing rows 54-89 down to 55-90.
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new daily record.
$ws.Range("A54").Value = 10
$ws.Range("B54").Value = "Vega Modelo de Temuco"
$ws.Range("C54").Value = "La Araucanía"
$ws.Range("D54").Value = 44777
$ws.Range("E54").Value = 9
$ws.Range("F54").Value = "Fruta"
$ws.Range("G54").Value = 100108
$ws.Range("H54").Value = "Tropicales y subtropicales"
$ws.Range("I54").Value = 100108007
$ws.Range("J54").Value = "Coco"
$ws.Range("K54").Value = "Sin especificar"
$ws.Range("L54").Value = "Primera"
$ws.Range("M54").Value = 30
$ws.Range("N54").Value = 30000
$ws.Range("O54").Value = 30000
$ws.Range("P54").Value = 30000
$ws.Range("Q54").Value = "$/malla 20 unidades"
$ws.Range("R54").Value = "Perú"
$ws.Range("S54").Value = 1500
$ws.Range("T54").Value = 20
